# GSKAU Wellness Main Shelf add in SOS and MSL
#
# Updates the "Facings SOS" include-list text (cells C2 and C3, which share
# the same value) to add "Wellness Main Shelf" and to use an en dash before
# "Grcy" for the NRT line. Also restores the active-cell selection on the
# frozen bottom-right pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dash = [char]0x2013
$newText = "Pain Main Shelf, Oral Main Shelf, Respiratory Main Shelf, NRT Main Shelf, Other Main Shelf, Pain Main Shelf - Grcy, Oral Main Shelf - Grcy, NRT Main Shelf " + $dash + " Grcy, Wellness Main Shelf"

$ws.Range("C2").Value = $newText
$ws.Range("C3").Value = $newText

# Row 3 re-wraps with the new text, shrinking its auto height.
$ws.Rows(3).RowHeight = 28.35

# Restore the selection on the bottom-right (frozen) pane to C4.
[void]$ws.Range("C4").Select()
